$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the "Condicion_Pacientes" table by one row and fill in the new
# daily COVID-19 patient condition data (2020-06-14 / serial 43996).
$lo = $ws.ListObjects.Item("Condicion_Pacientes")
$newRow = $lo.ListRows.Add()

$newRow.Range.Item(1, 1).Value = 43996
$newRow.Range.Item(1, 2).Value = 1187
$newRow.Range.Item(1, 3).Value = 403
$newRow.Range.Item(1, 4).Value = 436
$newRow.Range.Item(1, 5).Value = 252
$newRow.Range.Item(1, 6).Value = 52

# Match number formatting / style of the row above (date style + centered
# number style) by copying formats down, same as Excel does automatically
# when a table grows.
$ws.Range("A93:F93").Copy()
$ws.Range("A94:F94").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active selection to the newly added cell, matching the
# saved workbook state.
[void]$ws.Range("F94").Select()
